$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.748.20'
$ws.Range("E2").Value = '  +0.26%  '

$ws.Range("D3").Value = '3.010.62'
$ws.Range("E3").Value = '  +2.09%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '380.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.75'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.75%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.546'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.48%  '

$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.604'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.87%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.82'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.140'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0847'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.16%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.88'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.48%  '

$ws.Range("D14").Value = '3.482.95'
$ws.Range("E14").Value = '  +2.17%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.57'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.79%  '

$ws.Range("D16").Value = '3.014.13'
$ws.Range("E16").Value = '  +2.08%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.971'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.92%  '

$ws.Range("D18").Value = '51.758.07'
$ws.Range("E18").Value = '  +0.42%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.49'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.59%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.95%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.33%  '

$ws.Range("D22").Value = '0.0₃0965'
$ws.Range("E22").Value = '  +1.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '264.87'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.41%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.80'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.06%  '

$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.173'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.59%  '

$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.36'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +16.99%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.48'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.84%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '26.22'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.86%  '

$ws.Range("E30").Value = '  -0.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.106'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.64%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.93%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.68'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.11%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '51.42'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.15%  '

$ws.Range("E35").Value = '  -3.72%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0447'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.74%  '

$ws.Range("E37").Value = '  +0.12%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.13'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.56%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '17.53'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.74%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.69'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.98%  '

$ws.Range("E41").Value = '  -1.25%  '

$ws.Range("E42").Value = '  +1.93%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '124.49'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.39%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.49'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.69%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.282'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +18.58%  '

$ws.Range("E46").Value = '  -2.70%  '

$ws.Range("E47").Value = '  +7.37%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.31'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.87%  '

$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '2.051.62'
$ws.Range("E49").Value = '  -2.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0353'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +10.91%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.877'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.09%  '
